$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Replace("52.103.87", "52.044.93") | Out-Null
$ws.Range("E2").Replace("  -0.06%  ", "  -0.04%  ") | Out-Null
$ws.Range("D3").Replace("2.843.60", "2.841.03") | Out-Null
$ws.Range("E3").Replace("  +1.81%  ", "  +1.80%  ") | Out-Null
$ws.Range("E4").Replace("  +0.06%  ", "  +0.02%  ") | Out-Null
$ws.Range("D5").Replace("362.62", "361.54") | Out-Null
$ws.Range("E5").Replace("  +6.06%  ", "  +5.34%  ") | Out-Null
$ws.Range("D6").Replace("113.70", "113.55") | Out-Null
$ws.Range("E6").Replace("  -2.91%  ", "  -2.54%  ") | Out-Null
$ws.Range("D7").Replace("0.567", "0.568") | Out-Null
$ws.Range("E7").Replace("  +4.51%  ", "  +4.91%  ") | Out-Null
$ws.Range("E8").Replace("  +0.04%  ", "  -0.02%  ") | Out-Null
$ws.Range("D9").Replace("0.604", "0.603") | Out-Null
$ws.Range("E9").Replace("  +4.07%  ", "  +4.09%  ") | Out-Null
$ws.Range("D10").Replace("41.76", "41.65") | Out-Null
$ws.Range("E10").Replace("  -1.83%  ", "  -1.50%  ") | Out-Null
$ws.Range("D11").Replace("0.0861", "0.0860") | Out-Null
$ws.Range("E11").Replace("  -0.98%  ", "  -0.83%  ") | Out-Null
$ws.Range("B12").Replace("TRON", "Chainlink") | Out-Null
$ws.Range("C12").Replace("https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link") | Out-Null
$ws.Range("D12").Replace("0.131", "20.06") | Out-Null
$ws.Range("E12").Replace("  +1.21%  ", "  -0.36%  ") | Out-Null
$ws.Range("B13").Replace("Chainlink", "TRON") | Out-Null
$ws.Range("C13").Replace("https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx") | Out-Null
$ws.Range("D13").Replace("20.02", "0.131") | Out-Null
$ws.Range("E13").Replace("  -0.47%  ", "  +1.17%  ") | Out-Null
$ws.Range("D14").Replace("7.79", "7.77") | Out-Null
$ws.Range("E14").Replace("  +1.60%  ", "  +1.66%  ") | Out-Null
$ws.Range("D15").Replace("3.291.47", "3.289.25") | Out-Null
$ws.Range("E15").Replace("  +2.04%  ", "  +1.75%  ") | Out-Null
$ws.Range("D16").Replace("2.841.04", "2.840.56") | Out-Null
$ws.Range("E16").Replace("  +2.27%  ", "  +1.25%  ") | Out-Null
$ws.Range("E17").Replace("  +2.10%  ", "  +2.16%  ") | Out-Null
$ws.Range("D18").Replace("51.996.28", "51.924.80") | Out-Null
$ws.Range("E18").Replace("  -0.06%  ", "  -0.07%  ") | Out-Null
$ws.Range("D19").Replace("7.45", "7.44") | Out-Null
$ws.Range("E19").Replace("  +6.77%  ", "  +6.64%  ") | Out-Null
$ws.Range("E20").Replace("  -2.42%  ", "  -2.31%  ") | Out-Null
$ws.Range("D21").Replace("13.55", "13.53") | Out-Null
$ws.Range("E21").Replace("  +0.83%  ", "  +1.05%  ") | Out-Null
$ws.Range("D22").Replace("0.0₃0993", "0.0₃0992") | Out-Null
$ws.Range("E22").Replace("  +0.85%  ", "  +0.71%  ") | Out-Null
$ws.Range("D23").Replace("70.15", "70.20") | Out-Null
$ws.Range("E23").Replace("  -0.31%  ", "  -0.06%  ") | Out-Null
$ws.Range("D24").Replace("267.04", "267.11") | Out-Null
$ws.Range("E24").Replace("  -4.05%  ", "  -3.83%  ") | Out-Null
$ws.Range("D25").Replace("2.85", "2.83") | Out-Null
$ws.Range("E25").Replace("  +1.80%  ", "  +1.31%  ") | Out-Null
$ws.Range("D26").Replace("27.11", "27.13") | Out-Null
$ws.Range("E26").Replace("  +0.72%  ", "  +1.06%  ") | Out-Null
$ws.Range("E27").Replace("  -0.03%  ", "  +0.01%  ") | Out-Null
$ws.Range("D28").Replace("10.42", "10.41") | Out-Null
$ws.Range("E28").Replace("  +1.89%  ", "  +2.23%  ") | Out-Null
$ws.Range("E29").Replace("  +1.38%  ", "  +1.40%  ") | Out-Null
$ws.Range("D30").Replace("53.55", "53.79") | Out-Null
$ws.Range("E30").Replace("  +6.23%  ", "  +6.88%  ") | Out-Null
$ws.Range("E31").Replace("  -2.19%  ", "  -2.13%  ") | Out-Null
$ws.Range("D32").Replace("34.09", "34.06") | Out-Null
$ws.Range("E32").Replace("  -2.97%  ", "  -2.51%  ") | Out-Null
$ws.Range("D33").Replace("5.90", "5.88") | Out-Null
$ws.Range("E33").Replace("  +3.41%  ", "  +3.13%  ") | Out-Null
$ws.Range("D34").Replace("0.0447", "0.0446") | Out-Null
$ws.Range("E34").Replace("  +20.95%  ", "  +20.15%  ") | Out-Null
$ws.Range("D35").Replace("5.31", "5.33") | Out-Null
$ws.Range("E35").Replace("  +6.33%  ", "  +7.01%  ") | Out-Null
$ws.Range("D36").Replace("0.0839", "0.0843") | Out-Null
$ws.Range("E36").Replace("  +1.95%  ", "  +2.38%  ") | Out-Null
$ws.Range("E37").Replace("  +0.13%  ", "  -0.01%  ") | Out-Null
$ws.Range("E38").Replace("  -0.20%  ", "  +0.15%  ") | Out-Null
$ws.Range("E39").Replace("  -2.64%  ", "  -2.52%  ") | Out-Null
$ws.Range("E40").Replace("  -4.00%  ", "  -3.45%  ") | Out-Null
$ws.Range("D41").Replace("24.31", "24.19") | Out-Null
$ws.Range("E41").Replace("  +3.22%  ", "  +2.03%  ") | Out-Null
$ws.Range("E42").Replace("  +1.72%  ", "  +1.84%  ") | Out-Null
$ws.Range("B43").Replace("Stacks", "Monero") | Out-Null
$ws.Range("C43").Replace("https://coinranking.com/coin/mMPrMcB7+stacks-stx", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr") | Out-Null
$ws.Range("D43").Replace("2.57", "127.85") | Out-Null
$ws.Range("E43").Replace("  -6.49%  ", "  -0.18%  ") | Out-Null
$ws.Range("B44").Replace("Monero", "Stacks") | Out-Null
$ws.Range("C44").Replace("https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "https://coinranking.com/coin/mMPrMcB7+stacks-stx") | Out-Null
$ws.Range("D44").Replace("127.75", "2.56") | Out-Null
$ws.Range("E44").Replace("  -0.11%  ", "  -6.52%  ") | Out-Null
$ws.Range("E45").Replace("  -3.67%  ", "  -3.62%  ") | Out-Null
$ws.Range("D46").Replace("2.121.75", "2.121.07") | Out-Null
$ws.Range("E46").Replace("  +0.58%  ", "  +0.60%  ") | Out-Null
$ws.Range("D47").Replace("3.39", "3.38") | Out-Null
$ws.Range("E47").Replace("  +1.02%  ", "  +1.12%  ") | Out-Null
$ws.Range("E48").Replace("  +1.02%  ", "  +1.01%  ") | Out-Null
$ws.Range("D49").Replace("1.01", "0.999") | Out-Null
$ws.Range("E49").Replace("  +9.93%  ", "  +10.02%  ") | Out-Null
$ws.Range("E50").Replace("  +4.56%  ", "  +4.83%  ") | Out-Null
$ws.Range("D51").Replace("9.02", "9.03") | Out-Null
$ws.Range("E51").Replace("  +0.89%  ", "  +1.06%  ") | Out-Null
